# "Implemented the label and table tag"
#
# Adds two new locator sheets (LabelLocators, TableLocators) after the
# existing LinkLocators sheet, each seeded with the same 7-column header
# row (Loc1..Loc7) used by every other locator sheet in this workbook.
# Also nudges a couple of pre-existing selections, matching the author's
# commit.

$wb = $excel.ActiveWorkbook

$inputLocators = $wb.Worksheets.Item(1)
$linkLocators  = $wb.Worksheets.Item(5)

# --- New sheet: LabelLocators (inserted right after LinkLocators) ---
$labelLocators = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $linkLocators)
$labelLocators.Name = "LabelLocators"
$labelLocators.Range("A1").Value = "Loc1"
$labelLocators.Range("B1").Value = "Loc2"
$labelLocators.Range("C1").Value = "Loc3"
$labelLocators.Range("D1").Value = "Loc4"
$labelLocators.Range("E1").Value = "Loc5"
$labelLocators.Range("F1").Value = "Loc6"
$labelLocators.Range("G1").Value = "Loc7"

# --- New sheet: TableLocators (inserted right after LabelLocators) ---
$tableLocators = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $labelLocators)
$tableLocators.Name = "TableLocators"
$tableLocators.Range("A1").Value = "Loc1"
$tableLocators.Range("B1").Value = "Loc2"
$tableLocators.Range("C1").Value = "Loc3"
$tableLocators.Range("D1").Value = "Loc4"
$tableLocators.Range("E1").Value = "Loc5"
$tableLocators.Range("F1").Value = "Loc6"
$tableLocators.Range("G1").Value = "Loc7"

# --- Selection tweaks on existing sheets ---
$inputLocators.Range("B20").Select()
$linkLocators.Range("A1:G1").Select()
$labelLocators.Range("A1:G1").Select()

# TableLocators is selected last so it ends up the active tab/sheet,
# matching the workbook's new activeTab.
$tableLocators.Range("K25").Select()
